# Append two new log rows (Excel rows 158 and 159) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 158 -> run_id 157
$ws.Cells.Item(158, 1).Value = 157
$ws.Cells.Item(158, 2).Value = 1
$ws.Cells.Item(158, 3).Value = "2024-06-18 07:12:18"
$ws.Cells.Item(158, 4).Value = 200
$ws.Cells.Item(158, 5).Value = 6

# New row 159 -> run_id 158
$ws.Cells.Item(159, 1).Value = 158
$ws.Cells.Item(159, 2).Value = 2
$ws.Cells.Item(159, 3).Value = "2024-06-18 07:12:18"
$ws.Cells.Item(159, 4).Value = 200
$ws.Cells.Item(159, 5).Value = 0
